$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# Replace the "missing font" demo text with the "font substituted" message
# and switch the run's font to the (intentionally missing) "Beauty" typeface,
# mirroring how PowerPoint falls back to Courier New when a font isn't installed.
$tr.Text = "Courier New font is used instead of missing Beauty font"
$tr.Font.Name = "Beauty"
